$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(127).Insert()

$ws.Cells.Item(127,1).Value = 7
$ws.Cells.Item(127,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(127,3).Value = "Ñuble"
$ws.Cells.Item(127,4).Value = 44510
$ws.Cells.Item(127,5).Value = 16
$ws.Cells.Item(127,6).Value = 100112009
$ws.Cells.Item(127,7).Value = "Acelga"
$ws.Cells.Item(127,8).Value = "Sin especificar"
$ws.Cells.Item(127,9).Value = "Primera"
$ws.Cells.Item(127,10).Value = 100
$ws.Cells.Item(127,11).Value = 350
$ws.Cells.Item(127,12).Value = 400
$ws.Cells.Item(127,13).Value = 375
$ws.Cells.Item(127,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(127,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(127,16).Value = 375
$ws.Cells.Item(127,17).Value = 1
$ws.Cells.Item(127,18).Value = "Hortaliza"
